$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Metadata sheet: refresh the "Status" and "Date" property values
#    (IG re-deploy: status moves from "active" to "draft" and the
#    publication date is bumped to the new deploy timestamp).
# ------------------------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B6").Value = "draft"
$meta.Range("B8").Value = "2023-08-01T16:12:28+00:00"

# ------------------------------------------------------------------
# 2. Re-assert the vertical-top / wrap-text alignment on every cell
#    that carries it, so the stored style definitions explicitly flag
#    applyAlignment="true" (previously the alignment element was
#    present but not marked as applied) on every sheet that shares
#    those cell styles. Target the exact populated ranges on each
#    sheet so we don't materialize brand-new blank cells.
# ------------------------------------------------------------------
$xlVAlignTop = -4160

$meta.Range("A1:B14").WrapText = $true
$meta.Range("A1:B14").VerticalAlignment = $xlVAlignTop

$codes = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")
$codes.Range("A1:A2").WrapText = $true
$codes.Range("A1:A2").VerticalAlignment = $xlVAlignTop
$codes.Range("A3:B4").WrapText = $true
$codes.Range("A3:B4").VerticalAlignment = $xlVAlignTop
